$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-completed NPC row (row 5, "Cyborg Drone"):
# columns C-G (2D Texture, 3D Mesh, Animations, Anim Clips in Blend, Prefab) = 1
# column H (Testing) stays empty
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

# Move the active selection to G6, matching where the user last clicked
$ws.Range("G6").Select()
